# Applies corrected financial figures to the IFRS company_list sheet
# (rows 2-9, columns D:AJ) per the "error solve ifrs list" fix.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (2014/12  (IFRS연결))
$ws.Range("D2").Value = 51296
$ws.Range("E2").Value = 2134
$ws.Range("F2").Value = 2134
$ws.Range("G2").Value = 1469
$ws.Range("H2").Value = 4551
$ws.Range("I2").Value = 4556
$ws.Range("J2").Value = -6
$ws.Range("K2").Value = 95114
$ws.Range("L2").Value = 42261
$ws.Range("M2").Value = 52853
$ws.Range("N2").Value = 52831
$ws.Range("O2").Value = 22
$ws.Range("P2").Value = 135
$ws.Range("Q2").Value = 1431
$ws.Range("R2").Value = 544
$ws.Range("S2").Value = -456
$ws.Range("T2").Value = 1856
$ws.Range("U2").Value = -425
$ws.Range("V2").Value = 17599
$ws.Range("W2").Value = 4.16
$ws.Range("X2").Value = 8.869999999999999
$ws.Range("Y2").Value = 9.869999999999999
$ws.Range("Z2").Value = 5.04
$ws.Range("AA2").Value = 79.95999999999999
$ws.Range("AB2").Value = 22648.6
$ws.Range("AC2").Value = 3632
$ws.Range("AD2").Value = 43.51
$ws.Range("AE2").Value = 45557
$ws.Range("AF2").Value = 3.47
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 135000000

# Row 3 (2015/12  (IFRS연결))
$ws.Range("D3").Value = 133447
$ws.Range("E3").Value = 371
$ws.Range("F3").Value = 371
$ws.Range("G3").Value = 27757
$ws.Range("H3").Value = 26857
$ws.Range("I3").Value = 27468
$ws.Range("J3").Value = -612
$ws.Range("K3").Value = 423614
$ws.Range("L3").Value = 240475
$ws.Range("M3").Value = 183139
$ws.Range("N3").Value = 165017
$ws.Range("O3").Value = 18123
$ws.Range("P3").Value = 191
$ws.Range("Q3").Value = -623
$ws.Range("R3").Value = 23686
$ws.Range("S3").Value = -5274
$ws.Range("T3").Value = 2595
$ws.Range("U3").Value = -3218
$ws.Range("V3").Value = 77463
$ws.Range("W3").Value = 0.28
$ws.Range("X3").Value = 20.13
$ws.Range("Y3").Value = 25.22
$ws.Range("Z3").Value = 10.36
$ws.Range("AA3").Value = 131.31
$ws.Range("AB3").Value = 79713.47
$ws.Range("AC3").Value = 17857
$ws.Range("AD3").Value = 7.84
$ws.Range("AE3").Value = 98436
$ws.Range("AF3").Value = 1.42
$ws.Range("AG3").Value = 500
$ws.Range("AH3").Value = 0.36
$ws.Range("AI3").Value = 3.05
$ws.Range("AJ3").Value = 189690043

# Row 4 (2016/12  (IFRS연결))
$ws.Range("D4").Value = 281027
$ws.Range("E4").Value = 1395
$ws.Range("F4").Value = 1395
$ws.Range("G4").Value = 898
$ws.Range("H4").Value = 208
$ws.Range("I4").Value = 1074
$ws.Range("J4").Value = -865
$ws.Range("K4").Value = 444585
$ws.Range("L4").Value = 233528
$ws.Range("M4").Value = 211057
$ws.Range("N4").Value = 183016
$ws.Range("O4").Value = 28041
$ws.Range("P4").Value = 191
$ws.Range("Q4").Value = 13554
$ws.Range("R4").Value = -8121
$ws.Range("S4").Value = 1224
$ws.Range("T4").Value = 4847
$ws.Range("U4").Value = 8707
$ws.Range("V4").Value = 71305
$ws.Range("W4").Value = 0.5
$ws.Range("X4").Value = 0.07000000000000001
$ws.Range("Y4").Value = 0.62
$ws.Range("Z4").Value = 0.05
$ws.Range("AA4").Value = 110.65
$ws.Range("AB4").Value = 82053.37
$ws.Range("AC4").Value = 561
$ws.Range("AD4").Value = 223.65
$ws.Range("AE4").Value = 110964
$ws.Range("AF4").Value = 1.13
$ws.Range("AG4").Value = 550
$ws.Range("AH4").Value = 0.44
$ws.Range("AI4").Value = 84.56
$ws.Range("AJ4").Value = 189690043

# Row 5 (2017/12  (IFRS연결))
$ws.Range("D5").Value = 292790
$ws.Range("E5").Value = 8813
$ws.Range("F5").Value = 8813
$ws.Range("G5").Value = 8250
$ws.Range("H5").Value = 4811
$ws.Range("I5").Value = 6398
$ws.Range("J5").Value = -1587
$ws.Range("K5").Value = 490489
$ws.Range("L5").Value = 238977
$ws.Range("M5").Value = 251512
$ws.Range("N5").Value = 225568
$ws.Range("O5").Value = 25945
$ws.Range("P5").Value = 191
$ws.Range("Q5").Value = 13124
$ws.Range("R5").Value = 2799
$ws.Range("S5").Value = -11167
$ws.Range("T5").Value = 6780
$ws.Range("U5").Value = 6344
$ws.Range("V5").Value = 60342
$ws.Range("W5").Value = 3.01
$ws.Range("X5").Value = 1.64
$ws.Range("Y5").Value = 3.13
$ws.Range("Z5").Value = 1.03
$ws.Range("AA5").Value = 95.02
$ws.Range("AB5").Value = 84799.73
$ws.Range("AC5").Value = 3344
$ws.Range("AD5").Value = 37.68
$ws.Range("AE5").Value = 136764
$ws.Range("AF5").Value = 0.92
$ws.Range("AG5").Value = 2000
$ws.Range("AH5").Value = 1.59
$ws.Range("AI5").Value = 51.57
$ws.Range("AJ5").Value = 189690043

# Row 6 (2018/12  (IFRS연결))
$ws.Range("D6").Value = 311556
$ws.Range("E6").Value = 11039
$ws.Range("F6").Value = 11039
$ws.Range("G6").Value = 23827
$ws.Range("H6").Value = 17482
$ws.Range("I6").Value = 17128
$ws.Range("K6").Value = 424067
$ws.Range("L6").Value = 198577
$ws.Range("M6").Value = 225489
$ws.Range("N6").Value = 199466
$ws.Range("P6").Value = 191
$ws.Range("Q6").Value = 16784
$ws.Range("R6").Value = 1678
$ws.Range("S6").Value = -18847
$ws.Range("T6").Value = 5404
$ws.Range("U6").Value = 11380
$ws.Range("V6").Value = 45166
$ws.Range("W6").Value = 3.54
$ws.Range("X6").Value = 5.61
$ws.Range("Y6").Value = 8.06
$ws.Range("Z6").Value = 3.82
$ws.Range("AA6").Value = 88.06
$ws.Range("AB6").Value = 91044.17
$ws.Range("AC6").Value = 8953
$ws.Range("AD6").Value = 11.78
$ws.Range("AE6").Value = 120938
$ws.Range("AF6").Value = 0.87
$ws.Range("AG6").Value = 2000
$ws.Range("AH6").Value = 1.9
$ws.Range("AI6").Value = 19.26
$ws.Range("AJ6").Value = 189690043

# Row 7 (2019/12(E)  (IFRS연결))
$ws.Range("D7").Value = 308001
$ws.Range("E7").Value = 7964
$ws.Range("G7").Value = 13464
$ws.Range("H7").Value = 9286
$ws.Range("I7").Value = 9500
$ws.Range("K7").Value = 447808
$ws.Range("L7").Value = 198413
$ws.Range("M7").Value = 249395
$ws.Range("N7").Value = 223350
$ws.Range("P7").Value = 190
$ws.Range("Q7").Value = 10290
$ws.Range("R7").Value = -11638
$ws.Range("S7").Value = -3117
$ws.Range("T7").Value = 5162
$ws.Range("U7").Value = 5045
$ws.Range("W7").Value = 2.59
$ws.Range("X7").Value = 3.02
$ws.Range("Y7").Value = 4.49
$ws.Range("Z7").Value = 2.13
$ws.Range("AA7").Value = 79.56
$ws.Range("AC7").Value = 4966
$ws.Range("AD7").Value = 21.65
$ws.Range("AE7").Value = 135420
$ws.Range("AF7").Value = 0.79
$ws.Range("AG7").Value = 2018
$ws.Range("AH7").Value = 1.88
$ws.Range("AI7").Value = 40.28

# Row 8 (2020/12(E)  (IFRS연결))
$ws.Range("D8").Value = 309178
$ws.Range("E8").Value = 9970
$ws.Range("G8").Value = 16436
$ws.Range("H8").Value = 11881
$ws.Range("I8").Value = 11788
$ws.Range("K8").Value = 455911
$ws.Range("L8").Value = 197245
$ws.Range("M8").Value = 258667
$ws.Range("N8").Value = 235229
$ws.Range("P8").Value = 190
$ws.Range("Q8").Value = 16531
$ws.Range("R8").Value = -11349
$ws.Range("S8").Value = -4423
$ws.Range("T8").Value = 4661
$ws.Range("U8").Value = 9036
$ws.Range("W8").Value = 3.23
$ws.Range("X8").Value = 3.84
$ws.Range("Y8").Value = 5.14
$ws.Range("Z8").Value = 2.63
$ws.Range("AA8").Value = 76.25
$ws.Range("AC8").Value = 6161
$ws.Range("AD8").Value = 17.61
$ws.Range("AE8").Value = 142622
$ws.Range("AF8").Value = 0.76
$ws.Range("AG8").Value = 2333
$ws.Range("AH8").Value = 2.15
$ws.Range("AI8").Value = 37.54

# Row 9 (2021/12(E)  (IFRS연결))
$ws.Range("D9").Value = 317680
$ws.Range("E9").Value = 11306
$ws.Range("G9").Value = 18223
$ws.Range("H9").Value = 13236
$ws.Range("I9").Value = 12955
$ws.Range("K9").Value = 466018
$ws.Range("L9").Value = 198259
$ws.Range("M9").Value = 267759
$ws.Range("N9").Value = 244501
$ws.Range("P9").Value = 190
$ws.Range("Q9").Value = 15090
$ws.Range("R9").Value = -10529
$ws.Range("S9").Value = -4966
$ws.Range("T9").Value = 4628
$ws.Range("U9").Value = 7147
$ws.Range("W9").Value = 3.56
$ws.Range("X9").Value = 4.17
$ws.Range("Y9").Value = 5.4
$ws.Range("Z9").Value = 2.87
$ws.Range("AA9").Value = 74.04000000000001
$ws.Range("AC9").Value = 6772
$ws.Range("AD9").Value = 16.02
$ws.Range("AE9").Value = 148243
$ws.Range("AF9").Value = 0.73
$ws.Range("AG9").Value = 2396
$ws.Range("AH9").Value = 2.21
$ws.Range("AI9").Value = 35.08

Write-Output "Updated rows 2-9 (D:AJ) with corrected IFRS figures."